$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a randomly re-masked "missing data" sample. Two rows from
# the previous mask (originally row 26 "RM 232" and row 28 "SC 92") are
# dropped entirely in the new sample, shifting every row below them up.
# Deleting row 26 first shifts "SC 92" up to row 27, so delete that next.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- Re-apply the per-cell masking differences for the surviving rows ---
# (cells set to "" reproduce the empty/masked-out inline-string cells;
#  numeric assignments reproduce values that became unmasked.)

$ws.Range("E5").Value = ""        # RM 14   - newly masked
$ws.Range("F7").Value = ""        # RM 32   - newly masked

$ws.Range("E11").Value = -7.9     # RM 58   - newly revealed

$ws.Range("C19").Value = 13.2     # RM 125  - newly revealed
$ws.Range("E19").Value = ""       # RM 125  - newly masked

$ws.Range("C21").Value = ""       # RM 135  - newly masked

$ws.Range("C23").Value = 12.2     # RM 140  - newly revealed
$ws.Range("E23").Value = -7       # RM 140  - newly revealed

$ws.Range("F24").Value = 16.78    # RM 142a - newly revealed

$ws.Range("E25").Value = -7.1     # RM 145  - newly revealed

# Row 26 is now "SC 5"
$ws.Range("B26").Value = ""       # SC 5    - newly masked

# Row 27 is now "SC 101" (previously held "SC 92"'s data before deletion)
$ws.Range("B27").Value = -20.4    # SC 101  - newly revealed
$ws.Range("C27").Value = ""       # SC 101  - newly masked
$ws.Range("E27").Value = ""       # SC 101  - newly masked

# Row 28 is now "SC 105"
$ws.Range("F28").Value = 17.44    # SC 105  - newly revealed

# Row 29 is now "SC 119"
$ws.Range("B29").Value = ""       # SC 119  - newly masked
$ws.Range("E29").Value = ""       # SC 119  - newly masked

# Row 30 is now "SC 120"
$ws.Range("F30").Value = ""       # SC 120  - newly masked

# Row 32 is now "SC 193"
$ws.Range("F32").Value = ""       # SC 193  - newly masked

# Row 33 is now "SC 232"
$ws.Range("C33").Value = 10.4     # SC 232  - newly revealed
$ws.Range("E33").Value = -10.7    # SC 232  - newly revealed
